$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bot login entries (username/password pairs), appended below the
# existing rows in columns B (username) and D (password).
$ws.Range("B8").Value = "bobbybob8005"
$ws.Range("D8").Value = "fakepass8005"
$ws.Range("B9").Value = "larryjerry220"
$ws.Range("D9").Value = "parksandrec"

# Selection moves to F7 as recorded in the saved view state.
$ws.Range("F7").Select()
